# Generate Report for Handback
# Adds a new handback record (e670c125-69a2-42de-8745-e26ce230c0fc) as row 3
# on the Overview, zh-cn and de-de sheets, mirroring the existing row 2 layout.

$wb = $excel.ActiveWorkbook

$hyperFont = 15570276   # OLE BGR for RGB(0x64,0x95,0xED) == ARGB FF6495ED used by the workbook's HyperLink style
$dateFmt   = "yyyy-mm-dd HH:mm:ss"

function Set-PlainCell($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-DateCell($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
    $ws.Range($addr).NumberFormat = $dateFmt
}

function Add-HandbackHyperlink($ws, $addr, $url, $display) {
    $ws.Hyperlinks.Add($ws.Range($addr), $url, "", "", $display)
    # Re-assert the workbook's existing hyperlink look (underline + CornflowerBlue)
    # instead of the engine's default theme hyperlink style, so the new cell matches
    # the formatting already used for row 2's links.
    $ws.Range($addr).Font.Name = "Calibri"
    $ws.Range($addr).Font.Underline = $true
    $ws.Range($addr).Font.Color = $hyperFont
}

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) — row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-PlainCell $wsOverview "A3" "e670c125-69a2-42de-8745-e26ce230c0fc.md"
Set-PlainCell $wsOverview "B3" "e2e\e670c125-69a2-42de-8745-e26ce230c0fc.md"
Add-HandbackHyperlink $wsOverview "B3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2fa56f608d1cd574ad662fa34717f0292661a74f/e2e/e670c125-69a2-42de-8745-e26ce230c0fc.md" "e2e\e670c125-69a2-42de-8745-e26ce230c0fc.md"
Set-PlainCell $wsOverview "C3" ".md"
Set-PlainCell $wsOverview "E3" "Handed back: in sync with en-US"
Set-PlainCell $wsOverview "F3" "Handed back: in sync with en-US"
Set-DateCell  $wsOverview "G3" "2016-09-09 12:04:17"

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) — row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-PlainCell $wsZhCn "A3" "e670c125-69a2-42de-8745-e26ce230c0fc.md"
Add-HandbackHyperlink $wsZhCn "A3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2fa56f608d1cd574ad662fa34717f0292661a74f/e2e/e670c125-69a2-42de-8745-e26ce230c0fc.md" "e670c125-69a2-42de-8745-e26ce230c0fc.md"
Set-PlainCell $wsZhCn "B3" ".md"
Set-PlainCell $wsZhCn "C3" "Handed back: in sync with en-US"
Set-PlainCell $wsZhCn "D3" "e2e"
Set-PlainCell $wsZhCn "E3" "ht"
Set-PlainCell $wsZhCn "F3" "True"
Set-PlainCell $wsZhCn "G3" "e670c125-69a2-42de-8745-e26ce230c0fc.dda2c8b89d594cb9da67715e13487049e0998dc9.zh-cn.xlf"
Set-DateCell  $wsZhCn "H3" "2016-09-09 12:03:58"
Set-PlainCell $wsZhCn "I3" "e670c125-69a2-42de-8745-e26ce230c0fc.md"
Add-HandbackHyperlink $wsZhCn "I3" "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/fa94d933dd202efdac3e60f1071fe2b31ab4eaff/e2e/e670c125-69a2-42de-8745-e26ce230c0fc.md" "e670c125-69a2-42de-8745-e26ce230c0fc.md"
Set-PlainCell $wsZhCn "J3" "e670c125-69a2-42de-8745-e26ce230c0fc.dda2c8b89d594cb9da67715e13487049e0998dc9.zh-cn.xlf"
Set-DateCell  $wsZhCn "K3" "2016-09-09 12:04:50"
Set-PlainCell $wsZhCn "L3" ""
Set-PlainCell $wsZhCn "M3" "True"
Set-PlainCell $wsZhCn "N3" ""
Set-PlainCell $wsZhCn "O3" "False"
Set-PlainCell $wsZhCn "P3" ""

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) — row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-PlainCell $wsDeDe "A3" "e670c125-69a2-42de-8745-e26ce230c0fc.md"
Add-HandbackHyperlink $wsDeDe "A3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2fa56f608d1cd574ad662fa34717f0292661a74f/e2e/e670c125-69a2-42de-8745-e26ce230c0fc.md" "e670c125-69a2-42de-8745-e26ce230c0fc.md"
Set-PlainCell $wsDeDe "B3" ".md"
Set-PlainCell $wsDeDe "C3" "Handed back: in sync with en-US"
Set-PlainCell $wsDeDe "D3" "e2e"
Set-PlainCell $wsDeDe "E3" "ht"
Set-PlainCell $wsDeDe "F3" "True"
Set-PlainCell $wsDeDe "G3" "e670c125-69a2-42de-8745-e26ce230c0fc.dda2c8b89d594cb9da67715e13487049e0998dc9.de-de.xlf"
Set-DateCell  $wsDeDe "H3" "2016-09-09 12:04:17"
Set-PlainCell $wsDeDe "I3" "e670c125-69a2-42de-8745-e26ce230c0fc.md"
Add-HandbackHyperlink $wsDeDe "I3" "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/76b1e4a5a76514ebe1931de99cf98702f7e082d2/e2e/e670c125-69a2-42de-8745-e26ce230c0fc.md" "e670c125-69a2-42de-8745-e26ce230c0fc.md"
Set-PlainCell $wsDeDe "J3" "e670c125-69a2-42de-8745-e26ce230c0fc.dda2c8b89d594cb9da67715e13487049e0998dc9.de-de.xlf"
Set-DateCell  $wsDeDe "K3" "2016-09-09 12:05:17"
Set-PlainCell $wsDeDe "L3" ""
Set-PlainCell $wsDeDe "M3" "True"
Set-PlainCell $wsDeDe "N3" ""
Set-PlainCell $wsDeDe "O3" "False"
Set-PlainCell $wsDeDe "P3" ""

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

Write-Output "Handback report row added for e670c125-69a2-42de-8745-e26ce230c0fc"
